$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.30"
$ws.Range("D3").Value = "'23.75"
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "'3.559"
$ws.Range("E4").Value = "3LEOLEO"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "'5.233"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "'0.05863"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.465"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.328"
$ws.Range("E8").Value = "7GateTokenGT"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.8093"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "'0.8975"
$ws.Range("E10").Value = "9FTXTokenFTT"
$ws.Range("B11").Value = "One"
$ws.Range("C11").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D11").Value = "'0.01038"
$ws.Range("E11").Value = "10OneONEBestin24h"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1378"
$ws.Range("E12").Value = "11WazirXWRX"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.07259"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"
$ws.Range("B14").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C14").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D14").Value = "'0.03075"
$ws.Range("E14").Value = "13LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B15").Value = "BitrueCoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D15").Value = "'0.03056"
$ws.Range("E15").Value = "14BitrueCoinBTR"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "'0.09319"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.854"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("B18").Value = "BitForexToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D18").Value = "'0.001542"
$ws.Range("E18").Value = "17BitForexTokenBF"
$ws.Range("B19").Value = "CoinExToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D19").Value = "'0.04718"
$ws.Range("E19").Value = "18CoinExTokenCET"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "'0.006241"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.001264"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "HotbitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D22").Value = "'0.004581"
$ws.Range("E22").Value = "21HotbitTokenHTB"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.00008698"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("D24").Value = "'2.175"
$ws.Range("D25").Value = "'0.3198"
$ws.Range("D28").Value = "'0.0002339"
$ws.Range("D40").Value = "'0.03784"
$ws.Range("D41").Value = "'0.006328"
$ws.Range("D43").Value = "'0.002699"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").Value = "'0.007050"
$ws.Range("D45").Value = "'0.00005412"
$ws.Range("D47").Value = "'0.5399"
$ws.Range("D48").Value = "'0.02135"
